# Burndown - Sprint 3: add an "Other Tasks" catch-all row and update a
# couple of day-4 actuals, per the "Meeting & Input Support Stuff" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New catch-all task row 21: 7 estimated hours, 6 of them burned on Day 4.
$ws.Range("A21").Value = "Other Tasks - added for clarity"
$ws.Range("B21").Value = 7
$ws.Range("G21").Value = 6

# Extra Day 4 progress recorded against two existing tasks.
$ws.Range("G12").Value = 1
$ws.Range("G19").Value = 2

# Leave the selection where the author last left it.
[void]$ws.Range("H24").Select()
